# This script applies the "Second Correction" -> "First Update" edition
# rename to the Square One Standard Statement of Work document.
#
# It touches two paragraphs:
#   1. The "First Edition, Second Correction" title paragraph, which
#      becomes two runs: "First Edition, " + "First Update".
#   2. The "We agree to these terms..." paragraph, whose trailing runs
#      (which used to spell out "First" / " Edition, " / "Second" /
#      " Correction (...)") are reworked into: a combined lead-in run
#      ending in "First Edition, ", a "First Update" run, a run for the
#      literal text before the URL code, a run for the new "1u" code,
#      and a closing ")." run.
#
# Because the Word object model here does not expose a "split run"
# primitive, each target paragraph's content (excluding its paragraph
# mark) is replaced wholesale via Range.InsertXML with the exact set of
# runs/formatting desired, while preserving the paragraph's own <w:pPr>
# (which lives outside the replaced range).

$d = $word.ActiveDocument

function New-WordOpenXmlPackage([string]$bodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphRuns($paragraph, [string]$runsXml) {
    # Replace a paragraph's runs (but not its paragraph mark / pPr) with
    # the supplied raw <w:r>...</w:r> XML.
    $range = $d.Range($paragraph.Range.Start, $paragraph.Range.End - 1)
    $range.InsertXML((New-WordOpenXmlPackage("<w:p>" + $runsXml + "</w:p>")))
}

# Locate the two paragraphs by their (pre-edit) text so the script is not
# brittle to unrelated paragraph-count changes elsewhere in the document.
$titlePara = $null
$agreementPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "First Edition, Second Correction") {
        $titlePara = $p
    } elseif ($t.StartsWith("We agree to these terms, along with the Square One Standard Contractor Terms,")) {
        $agreementPara = $p
    }
}

if ($titlePara -eq $null) {
    throw "Could not find the 'First Edition, Second Correction' title paragraph"
}
if ($agreementPara -eq $null) {
    throw "Could not find the 'We agree to these terms...' paragraph"
}

# --- Paragraph: "First Edition, Second Correction" -> two runs.
$titleRuns = ''
$titleRuns += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/></w:rPr><w:t xml:space="preserve">First Edition, </w:t></w:r>'
$titleRuns += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:b/></w:rPr><w:t>First Update</w:t></w:r>'
Set-ParagraphRuns $titlePara $titleRuns

# --- Paragraph: "We agree to these terms, ..." -> five runs.
$agreementRuns = ''
$agreementRuns += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve">We agree to these terms, along with the Square One Standard Contractor Terms, First Edition, </w:t></w:r>'
$agreementRuns += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t>First Update</w:t></w:r>'
$agreementRuns += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve"> (https://squareoneforms.com/contractor/1e</w:t></w:r>'
$agreementRuns += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t>1u</w:t></w:r>'
$agreementRuns += '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t>).</w:t></w:r>'
Set-ParagraphRuns $agreementPara $agreementRuns

Write-Host "Title paragraph now reads:" $d.Paragraphs(2).Range.Text
Write-Host "Agreement paragraph now reads:" $d.Paragraphs(3).Range.Text
